$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(45891, 4535.36362203551, 3965.43917814788, 4368, 4904.608781, -1.38815262031797)
    3  = @(45892, 868.670076800243, 1551.97527604965, 1944, 1707.034027, 18.5974677603918)
    4  = @(45893, 773.003570643336, 1544.16414808315, 1944, 1602.771699, 17.913844851659)
    5  = @(45894, 5341.91604070659, 4248.95851628629, 1944, 5689.115541, 110.506584024154)
    6  = @(45895, 5341.91564027273, 4221.98295904657, 1944, 5689.115126, 109.382601865577)
    7  = @(45896, 5476.20428821302, 4257.14860152243, 1944, 5817.597141, 110.605893929559)
    8  = @(45897, 5476.20428821302, 4236.66826046088, 1944, 5817.597141, 109.752546385328)
    9  = @(45898, 5476.20428821302, 4332.1173398875, 1944, 5817.597141, 113.729591361437)
    10 = @(45899, 948.243530910975, 1628.64337215601, 1944, 1809.960285, 22.76500526021)
    11 = @(45900, 832.128236149678, 1568.8582723399, 1944, 1682.35776, 19.7953248412592)
    12 = @(45901, 5776.55841866516, 5007.58509614765, 2952, 6391.177173, 111.258493770103)
    13 = @(45902, 5784.11073458408, 4989.95486182018, 2952, 6398.793545, 110.526569676504)
    14 = @(45903, 5784.11073458408, 5023.79587282004, 2952, 6398.793545, 111.936611801499)
    15 = @(45904, 5812.79253171584, 5083.1238634848, 2952, 6445.485443, 115.15903228204)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}
